# New crime data collected - weekly CompStat update (114th Precinct)
# Updates volume/date headers and the weekly/28-day/YTD crime stat table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: force a cell to hold a given value as TEXT, re-using the display
# style already used elsewhere on the sheet for "text placeholder" cells
# (e.g. the ones that show "0" or "***.*" instead of a number). We do this by
# writing the value with a leading apostrophe (forces text) and then copying
# over just the number-format/style from a known-good donor cell so we don't
# end up with a stray quote-prefix style.
# ---------------------------------------------------------------------------
function Set-TextValue($cellRef, $donorRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Helper: force a cell to hold a given value as a NUMBER, copying number
# format/style from a known-good donor cell (used when a cell flips from a
# text placeholder to a real numeric value).
# ---------------------------------------------------------------------------
function Set-NumberValue($cellRef, $donorRef, $value) {
    $ws.Range($cellRef).Value = $value
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Header text: volume/number and report date range
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# ---------------------------------------------------------------------------
# Row 14 (Murder) - F14 flips from a number to the "0" text placeholder
# ---------------------------------------------------------------------------
Set-TextValue "F14" "C14" "0"

# ---------------------------------------------------------------------------
# Row 15 (Rape) - D15/E15 flip from numbers to text placeholders
# ---------------------------------------------------------------------------
Set-TextValue "D15" "C14" "0"
Set-TextValue "E15" "E14" "***.*"
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = 18.181818181818
$ws.Range("L15").Value = -7.142857142857
$ws.Range("M15").Value = 85.714285714285

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 10
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 87
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = 10.126582278481
$ws.Range("L16").Value = 31.818181818181
$ws.Range("M16").Value = -11.224489795918
$ws.Range("N16").Value = -86.124401913875

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 55
$ws.Range("H17").Value = -34.545454545454
$ws.Range("I17").Value = 147
$ws.Range("J17").Value = 173
$ws.Range("K17").Value = -15.028901734104
$ws.Range("L17").Value = 25.641025641025
$ws.Range("M17").Value = 81.481481481481
$ws.Range("N17").Value = -41.2

# ---------------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -35
$ws.Range("I18").Value = 72
$ws.Range("J18").Value = 89
$ws.Range("K18").Value = -19.101123595505
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -56.626506024096
$ws.Range("N18").Value = -89.007633587786

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 75
$ws.Range("H19").Value = 20.967741935483
$ws.Range("I19").Value = 305
$ws.Range("J19").Value = 276
$ws.Range("K19").Value = 10.507246376811
$ws.Range("L19").Value = 41.860465116279
$ws.Range("M19").Value = 39.269406392694
$ws.Range("N19").Value = -16.893732970027

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 500
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 46.666666666666
$ws.Range("I20").Value = 105
$ws.Range("J20").Value = 78
$ws.Range("K20").Value = 34.615384615384
$ws.Range("L20").Value = 98.113207547169
$ws.Range("M20").Value = 34.615384615384
$ws.Range("N20").Value = -88.782051282051

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = 18.918918918918
$ws.Range("F21").Value = 175
$ws.Range("G21").Value = 174
$ws.Range("H21").Value = 0.574712643678
$ws.Range("I21").Value = 732
$ws.Range("J21").Value = 708
$ws.Range("K21").Value = 3.389830508474
$ws.Range("L21").Value = 26.424870466321
$ws.Range("M21").Value = 12.098009188361
$ws.Range("N21").Value = -74.360770577933

# ---------------------------------------------------------------------------
# Row 22 (Transit) - C22 flips from the "0" text placeholder to a number
# ---------------------------------------------------------------------------
Set-NumberValue "C22" "C23" 2
$ws.Range("F22").Value = 3
$ws.Range("I22").Value = 11
$ws.Range("K22").Value = 120
$ws.Range("L22").Value = 10
$ws.Range("M22").Value = -8.333333333333

# ---------------------------------------------------------------------------
# Row 23 (Housing)
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 4
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 70
$ws.Range("J23").Value = 73
$ws.Range("K23").Value = -4.109589041095
$ws.Range("L23").Value = 29.629629629629
$ws.Range("M23").Value = 62.790697674418

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 46
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = -4.166666666666
$ws.Range("F24").Value = 170
$ws.Range("G24").Value = 198
$ws.Range("H24").Value = -14.141414141414
$ws.Range("I24").Value = 759
$ws.Range("J24").Value = 819
$ws.Range("K24").Value = -7.326007326007
$ws.Range("L24").Value = -1.811125485122
$ws.Range("M24").Value = 57.468879668049

# ---------------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 62
$ws.Range("G25").Value = 78
$ws.Range("H25").Value = -20.512820512820
$ws.Range("I25").Value = 258
$ws.Range("J25").Value = 311
$ws.Range("K25").Value = -17.041800643086
$ws.Range("L25").Value = 4.878048780487
$ws.Range("M25").Value = -6.859205776173

# ---------------------------------------------------------------------------
# Row 26 (UCR Rape*) - D26/E26 flip from numbers to text placeholders
# ---------------------------------------------------------------------------
Set-TextValue "D26" "C14" "0"
Set-TextValue "E26" "E14" "***.*"
$ws.Range("I26").Value = 16
$ws.Range("K26").Value = 6.666666666666
$ws.Range("L26").Value = -27.272727272727

# ---------------------------------------------------------------------------
# Row 27 (Other Sex Crimes) - D27/E27 flip from text placeholders to numbers
# ---------------------------------------------------------------------------
Set-NumberValue "D27" "D23" 1
Set-NumberValue "E27" "E23" 100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 30
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = 36.363636363636
$ws.Range("L27").Value = -6.25

# ---------------------------------------------------------------------------
# Row 28 (Shooting Vic.)
# ---------------------------------------------------------------------------
$ws.Range("F28").Value = 2
$ws.Range("M28").Value = 22.222222222222

# ---------------------------------------------------------------------------
# Row 29 (Shooting Inc.)
# ---------------------------------------------------------------------------
$ws.Range("F29").Value = 2
$ws.Range("M29").Value = 42.857142857142

# ---------------------------------------------------------------------------
# Row 30 (Hate Crimes) - F30 flips from a number to the "0" text placeholder
# ---------------------------------------------------------------------------
Set-TextValue "F30" "C14" "0"
